$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $addr, $val)
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

$ws.Range("D2").Value = "69.353.38"
$ws.Range("E2").Value = "  -2.62%  "
$ws.Range("D3").Value = "3.683.20"
$ws.Range("E3").Value = "  -3.45%  "
Set-TextCell $ws "D4" "0.999"
$ws.Range("E4").Value = "  -0.01%  "
Set-TextCell $ws "D5" "688.39"
$ws.Range("E5").Value = "  -2.43%  "
Set-TextCell $ws "D6" "162.47"
$ws.Range("E6").Value = "  -5.50%  "
$ws.Range("D7").Value = "3.682.42"
$ws.Range("E7").Value = "  -3.41%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -4.83%  "
$ws.Range("E10").Value = "  -8.69%  "
Set-TextCell $ws "D11" "7.36"
$ws.Range("E11").Value = "  -3.66%  "
$ws.Range("E12").Value = "  -3.71%  "
$ws.Range("E13").Value = "  -5.51%  "
Set-TextCell $ws "D14" "33.54"
$ws.Range("E14").Value = "  -6.92%  "
$ws.Range("D15").Value = "4.303.60"
$ws.Range("E15").Value = "  -3.48%  "
$ws.Range("D16").Value = "3.682.29"
$ws.Range("E16").Value = "  -2.04%  "
$ws.Range("D17").Value = "69.403.03"
$ws.Range("E17").Value = "  -2.52%  "
$ws.Range("E18").Value = "  -1.24%  "
Set-TextCell $ws "D19" "16.32"
$ws.Range("E19").Value = "  -6.93%  "
Set-TextCell $ws "D20" "6.62"
$ws.Range("E20").Value = "  -7.63%  "
Set-TextCell $ws "D21" "483.65"
$ws.Range("E21").Value = "  -6.73%  "
$ws.Range("E22").Value = "  -6.72%  "
$ws.Range("E23").Value = "  -8.02%  "
Set-TextCell $ws "D24" "80.24"
$ws.Range("E24").Value = "  -4.83%  "
$ws.Range("D25").Value = "3.827.96"
$ws.Range("E25").Value = "  -3.48%  "
$ws.Range("E26").Value = "  -9.75%  "
$ws.Range("E27").Value = "  +0.00%  "
Set-TextCell $ws "D28" "11.46"
$ws.Range("E28").Value = "  -4.84%  "
Set-TextCell $ws "D29" "9.51"
$ws.Range("E29").Value = "  -8.62%  "
$ws.Range("E30").Value = "  -10.71%  "
$ws.Range("E31").Value = "  -10.40%  "
Set-TextCell $ws "D32" "6.83"
$ws.Range("E32").Value = "  -7.61%  "
Set-TextCell $ws "D33" "2.08"
$ws.Range("E33").Value = "  -7.78%  "
$ws.Range("E34").Value = "  -6.95%  "
Set-TextCell $ws "D35" "1.00"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  -3.56%  "
$ws.Range("D37").Value = "3.648.33"
$ws.Range("E37").Value = "  -3.41%  "
$ws.Range("E38").Value = "  -7.46%  "
Set-TextCell $ws "D39" "6.40"
$ws.Range("E39").Value = "  +7.02%  "
Set-TextCell $ws "D40" "2.35"
$ws.Range("E40").Value = "  -2.05%  "
Set-TextCell $ws "D41" "0.0932"
$ws.Range("E41").Value = "  -8.01%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("E43").Value = "  -0.06%  "
Set-TextCell $ws "D44" "0.951"
$ws.Range("E44").Value = "  -7.24%  "
Set-TextCell $ws "D45" "163.91"
$ws.Range("E45").Value = "  -4.21%  "
Set-TextCell $ws "D46" "47.87"
$ws.Range("E46").Value = "  -3.43%  "
$ws.Range("E47").Value = "  -13.66%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell $ws "D48" "29.81"
$ws.Range("E48").Value = "  +2.60%  "
$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
Set-TextCell $ws "D49" "0.000288"
$ws.Range("E49").Value = "  -8.13%  "
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextCell $ws "D50" "1.35"
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("E51").Value = "  -1.93%  "
